# Swap the order of "System" and the email address in the
# "Recorded By" column (column G) of the Session Analysis Results sheet.
# Every cell that currently reads "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$updated = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Host "Updated $updated cell(s) in column G."
